# Add new "mail" API routes to the Api Paths sheet (Sheet1)
# Values are populated column-by-column so that new entries land in the
# shared string table in the same order as the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = "/mail/send"
$ws.Range("A9").Value = "/mail/send-multiple"

$ws.Range("B8").Value = "post"
$ws.Range("B9").Value = "post"

$ws.Range("C8").Value = "send mail to a single user"
$ws.Range("C9").Value = "send mail to multiple users"

$ws.Range("D8").Value = "yes"
$ws.Range("D9").Value = "yes"

$ws.Range("E8").Value = "Email`nBody"
$ws.Range("E9").Value = "EmailAddresses`nBody"

$ws.Range("F8").Value = "status and message"
$ws.Range("F9").Value = "status and message"

# Match the existing "Body" style used for E4 (wrap text) on the new cells
$ws.Range("E8:E9").WrapText = $true

# Match the row height used by similar wrapped rows
$ws.Range("A8:F9").RowHeight = 30

# Update selection to reflect the last edited cell, as in the source workbook
$ws.Range("D8").Select()
